$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "18-08-2021"
$ws.Range("B5").Value = 6000
$ws.Range("C5").Value = 13490
$ws.Range("D5").Value = 6000
$ws.Range("E5").Value = 510
$ws.Range("F5").Value = 5490
$ws.Range("G5").Value = 1.39
